$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.890.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.45%  "

$ws.Range("D3").Value = "'3.148.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.67%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'572.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.26%  "

$ws.Range("D6").Value = "'149.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.10%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "'3.148.41"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.82%  "

$ws.Range("D9").Value = "'0.525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.47%  "

$ws.Range("E10").Value = "  +3.34%  "

$ws.Range("D11").Value = "'6.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").Value = "'0.496"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.75%  "

$ws.Range("D13").Value = "'0.0000263"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +12.94%  "

$ws.Range("D14").Value = "'36.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.02%  "

$ws.Range("D15").Value = "'3.667.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.96%  "

$ws.Range("D16").Value = "'64.983.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.60%  "

$ws.Range("D17").Value = "'3.146.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.74%  "

$ws.Range("D18").Value = "'7.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.10%  "

$ws.Range("D20").Value = "'504.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.11%  "

$ws.Range("D21").Value = "'14.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.82%  "

$ws.Range("D22").Value = "'0.714"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.10%  "

$ws.Range("D23").Value = "'15.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.43%  "

$ws.Range("D24").Value = "'7.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.08%  "

$ws.Range("D25").Value = "'84.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.68%  "

$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("E27").Value = "  +3.09%  "

$ws.Range("D28").Value = "'8.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.03%  "

$ws.Range("D29").Value = "'2.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.62%  "

$ws.Range("E30").Value = "  +7.66%  "

$ws.Range("D31").Value = "'27.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.68%  "

$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("E33").Value = "  +2.54%  "

$ws.Range("D34").Value = "'6.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.24%  "

$ws.Range("D35").Value = "'6.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.30%  "

$ws.Range("D36").Value = "'54.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("D37").Value = "'0.0895"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.51%  "

$ws.Range("D38").Value = "'463.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.82%  "

$ws.Range("E39").Value = "  +1.49%  "

$ws.Range("E40").Value = "  +7.65%  "

$ws.Range("D41").Value = "'8.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.11%  "

$ws.Range("D42").Value = "'3.048.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").Value = "'0.116"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("E44").Value = "  +7.06%  "

$ws.Range("E45").Value = "  +1.37%  "

$ws.Range("D46").Value = "'28.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.17%  "

$ws.Range("D47").Value = "'0.0₃0580"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +11.47%  "

$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("E50").Value = "  +3.70%  "

$ws.Range("D51").Value = "'119.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.40%  "
